# "Rewrite the create page": shift the 38-day date column (A2:A39) from
# starting 2000-01-01 to starting 2020-01-01, and update the sheet's
# active selection to the refreshed date range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New starting date serial (2020-01-01) -- old range started at 2000-01-01
# (serial 36526) and ran for 38 consecutive days; keep the same length,
# just offset every date by +7305 days (20 years later, same month/day).
$startSerial = 43831

for ($i = 0; $i -lt 38; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $startSerial + $i
}

# Match the author's final selection: A2:A39 active at A2.
[void]$ws.Range("A2:A39").Select()
